$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.42166531085968
$ws.Range("B1").Value = 1.979551792144775
$ws.Range("C1").Value = 2.98771595954895
$ws.Range("D1").Value = 4.882411479949951
$ws.Range("E1").Value = 0.9247080683708191
